# Update countries & provincias Spain
# - Refresh the COVID stats snapshot (timestamp + several country rows).
# - A handful of neighbouring rows swap which country occupies which row
#   because the new snapshot re-sorted the ranking (the country that moved
#   up takes the higher row with fresh numbers; the other keeps its old,
#   unchanged numbers one row down).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Timestamp banner (row 1)
$ws.Range("A1").Value = "Datos actualizados a 12 de Septiembre de 2020 a las 11:46"

function Set-Row([int]$r, [string]$country, $total, $nuevos, $activos, $recuperados, $criticos, $muertesHoy, $muertes) {
    $ws.Range("A$r").Value = $country
    $ws.Range("B$r").Value = $total
    $ws.Range("C$r").Value = $nuevos
    $ws.Range("D$r").Value = $activos
    $ws.Range("E$r").Value = $recuperados
    $ws.Range("F$r").Value = $criticos
    $ws.Range("G$r").Value = $muertesHoy
    $ws.Range("H$r").Value = $muertes
}

# Straight data refreshes (country stays put, numbers updated)
Set-Row 4  "Estados Unidos" 6636580 333  3918491 2520665 0 3   197424
Set-Row 18 "Banglades"      336044  1282 238271  93071   0 34  4702
Set-Row 26 "Indonesia"      214746  3806 152458  53638   0 106 8650
Set-Row 49 "Polonia"        73650   603  59725   11743   0 13  2182
Set-Row 61 "Suiza"          46239   0    38500   5719    0 0   2020
Set-Row 95 "Consejo Danes para los Refugiados" 10385 24 9719 404 0 0 262
Set-Row 97 "Malasia"        9868    58   9189    551     0 0   128

# Rows that swap ranking position with their neighbour, each gaining a
# refreshed snapshot while the displaced country keeps its old figures.
Set-Row 70  "Austria"             32696 869 26579 5363 0 4 754
Set-Row 71  "Serbia"              32228 0   31100 398  0 0 730

Set-Row 111 "Eslovaquia"          5453  201 3114  2301 0 1 38
Set-Row 112 "Republica de Yibuti" 5394  0   5327  6    0 0 61

Set-Row 125 "Eslovenia"           3603  105 2699  769  0 0 135
Set-Row 126 "Jamaica"             3511  0   1019  2452 0 0 40

Set-Row 132 "Lituania"            3296  53  2070  1140 0 0 86
Set-Row 133 "Angola"              3279  0   1288  1860 0 0 131
